$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 date serial moved forward one day (2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Updated prices
$ws.Range("D33").Value = 954.2
$ws.Range("D34").Value = 1200
$ws.Range("D35").Value = 1231.1
$ws.Range("D36").Value = 1378
